$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "24.440.99"
$ws.Range("E2").Value = "  -1.89%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.655.80"
$ws.Range("E3").Value = "  -3.33%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.16%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.31"
$ws.Range("E5").Value = "  +0.63%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.000"
$ws.Range("E6").Value = "  +0.25%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3636"
$ws.Range("E7").Value = "  -3.13%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "47.14"
$ws.Range("E8").Value = "  -4.84%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3276"
$ws.Range("E9").Value = "  -5.18%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.128"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07076"
$ws.Range("E11").Value = "  -6.39%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.19%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "5.953"
$ws.Range("E13").Value = "  -5.82%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "19.50"
$ws.Range("E14").Value = "  -7.78%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.619"
$ws.Range("E15").Value = "  -6.24%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.656.16"
$ws.Range("E16").Value = "  -3.34%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001050"
$ws.Range("E17").Value = "  -7.71%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06605"
$ws.Range("E18").Value = "  -1.88%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9997"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "78.30"
$ws.Range("E20").Value = "  -7.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.946"
$ws.Range("E21").Value = "  -7.12%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.75"
$ws.Range("E22").Value = "  -9.38%  "

$ws.Range("E23").Value = "  -5.62%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "24.439.71"
$ws.Range("E24").Value = "  -1.83%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.478"
$ws.Range("E25").Value = "  +1.14%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.356"
$ws.Range("E26").Value = "  -15.90%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.35"
$ws.Range("E27").Value = "  -3.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.65"
$ws.Range("E28").Value = "  -8.80%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.841.79"
$ws.Range("E29").Value = "  -3.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.076"
$ws.Range("E32").Value = "  -4.57%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.722"
$ws.Range("E33").Value = "  -17.89%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08470"
$ws.Range("E34").Value = "  -4.20%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.668"
$ws.Range("E35").Value = "  -9.97%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "12.39"
$ws.Range("E36").Value = "  -11.08%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.238"
$ws.Range("E37").Value = "  -7.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06080"
$ws.Range("E38").Value = "  -9.30%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02218"
$ws.Range("E39").Value = "  -8.55%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2079"
$ws.Range("E40").Value = "  -7.71%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "8.251"
$ws.Range("E41").Value = "  -10.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.203"
$ws.Range("E42").Value = "  -6.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9996"
$ws.Range("E43").Value = "  +0.20%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5946"
$ws.Range("E44").Value = "  -8.36%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.796"
$ws.Range("E45").Value = "  -0.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.66"
$ws.Range("E46").Value = "  -9.57%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5656"
$ws.Range("E47").Value = "  -8.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.74"
$ws.Range("E48").Value = "  -6.00%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.952"
$ws.Range("E49").Value = "  -9.21%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06914"
$ws.Range("E50").Value = "  -5.76%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.97"
$ws.Range("E51").Value = "  -6.50%  "

# Swap rows 30 and 31: ImmutableX now at 30, BitcoinCash now at 31
$ws.Range("B30").Value = "ImmutableX"
$ws.Range("C30").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.205"
$ws.Range("E30").Value = "  -3.46%  "

$ws.Range("B31").Value = "BitcoinCash"
$ws.Range("C31").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "124.62"
$ws.Range("E31").Value = "  -6.09%  "
